# fixed issue with L1 on BoM
# L1 (row 16) is an Inductor whose part was changed from the TDK
# MLG0603P9N1HT000 to the Murata LQW18AN9N1D00D, along with an updated
# unit price and Mouser product link.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A16").Value = "LQW18AN9N1D00D"
$ws.Range("H16").Value = 0.22
$ws.Range("J16").Value = "https://au.mouser.com/ProductDetail/Murata-Electronics/LQW18AN9N1D00D?qs=EPZCdHdMnYJMYfHLXCR5Dg%3D%3D&_gl=1*o1ihxe*_ga*dW5kZWZpbmVk*_ga_15W4STQT4T*dW5kZWZpbmVk*_ga_1KQLCYKRX3*dW5kZWZpbmVk"

# Update the active selection to A16, matching the saved view state.
$ws.Range("A16").Select()
